# Finalization of 10 years of data:
# The daily weather table (A9:K40) on the original "Data Harian - Table"
# sheet gets selected/copied and pasted (as a static, formatted snapshot)
# onto a brand-new "Sheet1", which becomes the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Keep gridlines showing (engine quirk: the source file's showGridLines="true"
# literal is otherwise mis-read as off after a resave) and mirror the
# original author's on-screen selection of the table before copying it.
$excel.ActiveWindow.DisplayGridlines = $true
$ws1.Range("A9:K40").Select()

# Add the destination sheet right after the existing one, so tab order
# becomes: "Data Harian - Table", "Sheet1".
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Copy-and-paste the table (values + formatting) onto the new sheet,
# anchored at A1 so it becomes a standalone A1:K32 block.
$ws1.Range("A9:K40").Copy($ws2.Range("A1"))

# Match the row height Excel computes once the table no longer has the
# wide (13-char) source columns, which makes the wrapped header-style text
# wrap onto a second line for every data row.
$ws2.Rows("2:32").RowHeight = 28.8

# The newly pasted range ends up selected, and the new sheet becomes the
# active / visible tab.
$ws2.Range("A1:K32").Select()
$ws2.Activate()
